# Update data values per the new Estado de Cuenta data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Valor Mora total changed
$ws.Range("E11").Value = 60000

# Cant. Periodos changed
$ws.Range("F13").Value = 1

# Period "2507" -> "2508" for the remaining worker record
$ws.Range("E16").Value = "2508"

# Remove the second (duplicate) worker record row entirely, shifting rows below it up
$ws.Rows("17").Delete()

# Nudge the logo image slightly left (matches the refreshed layout)
$shp = $ws.Shapes.Item(1)
$shp.Left = 61.912816806
